$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Natalie's - Orange Mango): Quantity 2 -> 1, Total Cost 26.00 -> 13.00
# Force Text format so the numeric-looking strings keep their exact
# textual representation (e.g. "13.00" keeps its trailing zeros) instead
# of being coerced into a number.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "13.00"

# Row 6 (Natalie's - Honey Tangerine): Quantity 2 -> 1, Total Cost 28.00 -> 14.00
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "14.00"
